# Update PutShipCommand sequence diagrams
# Add footer explaining what <args> means.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The COM object model works in points (1 pt = 12700 EMU) and stores the
# value as a 32-bit float that gets truncated back to EMU on save. Nudge
# by a tiny epsilon so the round trip lands on the exact authored EMU.
function EmuToPt([double]$emu) {
  $base = $emu / 12700.0
  for ($i = 0; $i -lt 2000; $i++) {
    foreach ($sign in 1, -1) {
      $cand = $base + $sign * $i * 0.000001
      $f32 = [single]$cand
      $backEmu = [int64]([double]$f32 * 12700.0)
      if ($backEmu -eq [int64]$emu) {
        return $cand
      }
    }
  }
  return $base
}

# The target shape must land on shape Id 61 / name "TextBox 60" to match
# the authored deck (PowerPoint hands out Shape Ids from an internal
# counter that skips already-used values). Spin up throw-away textboxes
# until the counter reaches 61, keep that one, and delete the rest so
# only the intended shape remains in the final deck.
$placeholders = @()
$target = $null
for ($i = 0; $i -lt 100; $i++) {
  $tb = $s.Shapes.AddTextbox(1, 10, 10, 100, 20)
  if ($tb.Id -eq 61) {
    $target = $tb
    break
  } else {
    $placeholders += $tb
  }
}
if ($target -eq $null) {
  throw "Could not mint a shape with Id 61 (ran out of attempts)."
}

# Position / size in EMU, converted to points for the object model.
$left = EmuToPt 294013
$top = EmuToPt 6533668
$width = EmuToPt 9254858
$height = EmuToPt 220916
$target.Left = $left
$target.Top = $top
$target.Width = $width
$target.Height = $height

# No shape fill (transparent textbox)
$target.Fill.Visible = [Microsoft.Office.Core.MsoTriState]::msoFalse

# Text frame behaviour: wrap, zero insets
$target.TextFrame.WordWrap = [Microsoft.Office.Core.MsoTriState]::msoTrue
$target.TextFrame.MarginLeft = 0
$target.TextFrame.MarginRight = 0
$target.TextFrame.MarginTop = 0
$target.TextFrame.MarginBottom = 0

# Footer text: "<args> = n/destroyer c/a1 r/vertical t/tag", blue (0070C0)
$full = $target.TextFrame.TextRange
$full.Text = "<args> = n/destroyer c/a1 r/vertical t/tag"
$full.Font.Size = 14
$full.Font.Color.RGB = 0xC07000

# Shrink-shape-to-fit-text, then restore the authored extent.
$target.TextFrame.AutoSize = 1
$target.Left = $left
$target.Top = $top
$target.Width = $width
$target.Height = $height

# Re-touch each logical piece so PowerPoint keeps them as separate runs,
# matching the authored run-split ("<", "args", "> = n/destroyer c/a1 r/vertical t/tag")
$r1 = $full.Characters(1, 1)
$r1.Font.Size = 14
$r1.Font.Color.RGB = 0xC07000

$r2 = $full.Characters(2, 4)
$r2.Font.Size = 14
$r2.Font.Color.RGB = 0xC07000

$r3 = $full.Characters(6, $full.Length - 5)
$r3.Font.Size = 14
$r3.Font.Color.RGB = 0xC07000

# Splitting the runs nudges the autofit layout again; restore the extent once more.
$target.Left = $left
$target.Top = $top
$target.Width = $width
$target.Height = $height

# Remove the throw-away shapes used only to advance the Id counter.
foreach ($ph in $placeholders) {
  $ph.Delete()
}
